{"js": "// UC014 - Minha Conta Banc\u00e1ria: bump version history row + small wording fixes.\nconst replacements = [\n  // Version history table (row 1.0 -> 1.2.5)\n  [\"1.0\", \"1.2.5\"],\n  [\"Creation\", \"Update\"],\n  [\"Fabr\u00edcio Ara\u00fajo\", \"Julio Paiva\"],\n  [\"09/07/2020\", \"31/05/2023\"],\n  // Preconditions cell: fix typo \"usuario\" -> \"usu\u00e1rio\" and add trailing period\n  [\n    \"O usuario devidamente autenticado e na tela inicial do sistema\",\n    \"O usu\u00e1rio devidamente autenticado e na tela inicial do sistema.\",\n  ],\n  // Main flow step 1: add trailing period\n  [\n    \"1. Chefe Acessa a funcionalidade Minha Conta Banc\u00e1ria (menu) \",\n    \"1. Chefe Acessa a funcionalidade Minha Conta Banc\u00e1ria (menu). \",\n  ],\n  // Main flow step 4: fix \"banc\u00e1rios\" -> \"banc\u00e1ria\" (agreement with \"conta\")\n  [\n    \"4. System Exibe mensagens informativas (MSG403 - Informativos sobre a atualiza\u00e7\u00e3o de conta banc\u00e1rios (dados banc\u00e1rios)) para o usu\u00e1rio sobre a manuten\u00e7\u00e3o de informa\u00e7\u00f5es banc\u00e1rias. \",\n    \"4. System Exibe mensagens informativas (MSG403 - Informativos sobre a atualiza\u00e7\u00e3o de conta banc\u00e1ria (dados banc\u00e1rios)) para o usu\u00e1rio sobre a manuten\u00e7\u00e3o de informa\u00e7\u00f5es banc\u00e1rias. \",\n  ],\n  // Alternate flow AF[1] step 2: add trailing period\n  [\n    \"2. System Apresenta os campos (banco/ag\u00eancia/conta corrente) alterados \",\n    \"2. System Apresenta os campos (banco/ag\u00eancia/conta corrente) alterados. \",\n  ],\n];\n\nfor (const [searchText, replaceText] of replacements) {\n  const results = context.document.body.search(searchText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: ${JSON.stringify(searchText)}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# UC014 - Minha Conta Banc\u00e1ria: bump version history row + small wording fixes.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    # Version history table (row 1.0 -> 1.2.5)\n    @{Find = \"1.0\"; Replace = \"1.2.5\"},\n    @{Find = \"Creation\"; Replace = \"Update\"},\n    @{Find = \"Fabr\u00edcio Ara\u00fajo\"; Replace = \"Julio Paiva\"},\n    @{Find = \"09/07/2020\"; Replace = \"31/05/2023\"},\n    # Preconditions cell: fix typo \"usuario\" -> \"usu\u00e1rio\" and add trailing period\n    @{Find = \"O usuario devidamente autenticado e na tela inicial do sistema\"; Replace = \"O usu\u00e1rio devidamente autenticado e na tela inicial do sistema.\"},\n    # Main flow step 1: add trailing period\n    @{Find = \"1. Chefe Acessa a funcionalidade Minha Conta Banc\u00e1ria (menu) \"; Replace = \"1. Chefe Acessa a funcionalidade Minha Conta Banc\u00e1ria (menu). \"},\n    # Main flow step 4: fix \"banc\u00e1rios\" -> \"banc\u00e1ria\" (agreement with \"conta\")\n    @{Find = \"4. System Exibe mensagens informativas (MSG403 - Informativos sobre a atualiza\u00e7\u00e3o de conta banc\u00e1rios (dados banc\u00e1rios)) para o usu\u00e1rio sobre a manuten\u00e7\u00e3o de informa\u00e7\u00f5es banc\u00e1rias. \"; Replace = \"4. System Exibe mensagens informativas (MSG403 - Informativos sobre a atualiza\u00e7\u00e3o de conta banc\u00e1ria (dados banc\u00e1rios)) para o usu\u00e1rio sobre a manuten\u00e7\u00e3o de informa\u00e7\u00f5es banc\u00e1rias. \"},\n    # Alternate flow AF[1] step 2: add trailing period\n    @{Find = \"2. System Apresenta os campos (banco/ag\u00eancia/conta corrente) alterados \"; Replace = \"2. System Apresenta os campos (banco/ag\u00eancia/conta corrente) alterados. \"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.Replace\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $found = $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $true, 1, $false, $r.Replace, 2)\n    if (-not $found) {\n        throw \"Find/Replace failed for: $($r.Find)\"\n    }\n}\n"}
